$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MotorWare USER Variables")

# B13: switch to RC pole-based formula
$ws.Range("B13").Formula = "=1/(2*PI()*10*10^-9*1/(1/30000+1/2000000))"

# B15: switch from a static value to a formula
$ws.Range("B15").Formula = "=(30+2000)/30*(4.98*3.3/3.18)"

# New helper formulas in row 27-29, column F
$ws.Range("F27").Formula = "=(30+2000)/30"
$ws.Range("F28").Formula = "=B15/F27"
$ws.Range("F29").Formula = "=F28*3.18/4.98"

# Move the active selection
$ws.Range("B27").Select()

$wb.Save()
